$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    5373.441499999999,
    5250.357,
    5171.886999999999,
    5107.980499999999,
    5104.393,
    5230.476999999999,
    5902.879499999999,
    7037.001999999999,
    8592.440500000001,
    14095.0215,
    15828.8445,
    15449.287,
    15803.2595,
    15862.028,
    16032.3205,
    16559.4275,
    16683.57249999999,
    16252.82049999999,
    15546.958,
    13748.4165,
    12030.7215,
    9319.645999999999,
    6308.333499999995,
    5687.261999999999
)

$row = 2
foreach ($val in $values) {
    $ws.Cells.Item($row, 2).Value = $val
    $row++
}
